$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" placeholder text (15.09.2020 ->
#    16.09.2020) wherever it appears: the slide master and every slide
#    layout's date placeholder.
# ---------------------------------------------------------------------
$oldDate = "15.09.2020"
$newDate = "16.09.2020"

$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Fix the "6 ... Planned Motions" agenda entry on the last slide:
#    " Online " + "Stabilization" + " " + "of" + " the " + "Planned" + " " + "Motions"
#    ->
#    " Validation " + "of" + " " + "Planned" + " " + "Motions"
#    Edits are applied right-to-left so earlier character offsets stay
#    valid.
# ---------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)
$agendaShape = $null
for ($j = 1; $j -le $lastSlide.Shapes.Count; $j++) {
    $sh = $lastSlide.Shapes.Item($j)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "6 Online Stabilization of the Planned Motions") {
            $agendaShape = $sh
        }
    }
}

if ($agendaShape -ne $null) {
    $tr = $agendaShape.TextFrame.TextRange
    # remove "of the " (chars 24..30)
    $tr.Characters(24, 7).Text = ""
    # "Stabilization" (chars 10..22) -> "of"
    $tr.Characters(10, 13).Text = "of"
    # " Online " (chars 2..9) -> " Validation "
    $tr.Characters(2, 8).Text = " Validation "
}

# ---------------------------------------------------------------------
# 3) Add the new 9th slide: a duplicate of the (now corrected) last
#    slide, appended at the end of the deck.
# ---------------------------------------------------------------------
$lastSlide.Duplicate() | Out-Null
